$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Update the first three rows to "0M"
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# 2. Insert 10 new rows before the (old) 4th row, which will become rows 4-13
$newRowCount = 10
for ($i = 0; $i -lt $newRowCount; $i++) {
    $t.Rows.Add($t.Rows.Item(4)) | Out-Null
}

# 3. Populate the newly inserted rows with their values
$newVals = @("703", "0.00002", "0.00006", "0.00004", "0.00000", "0.00004", "0.00004", "0.00004", "0.02570", "100.0")
for ($i = 0; $i -lt $newVals.Length; $i++) {
    $t.Cell(4 + $i, 1).Range.Text = $newVals[$i]
}

# 4. Simplify the three rows that used to hold tab-separated multi-run values.
#    Those rows were originally rows 34, 35, 36 (1-indexed); after inserting
#    10 rows above them they are now rows 44, 45, 46.
$t.Cell(44,1).Range.Text = "99.99"
$t.Cell(45,1).Range.Text = "0.03"
$t.Cell(46,1).Range.Text = "333"
